$wb = $excel.ActiveWorkbook

# The "Feedback" sheet has a "Name" column (column B) that records each
# participant's name. Remove it entirely (shifting the remaining columns
# left) so individual respondents can no longer be identified.
$ws = $wb.Worksheets.Item("Feedback")
$ws.Activate()
$ws.Columns.Item(2).Delete() | Out-Null

# Leave the active selection/cursor where it ended up after the edit.
$ws.Range("D11").Select() | Out-Null
